# Update countries & provincias Spain
#
# The "Pais" sheet lists COVID-19 stats per country. Column A cells hold
# the country name (looked up from the shared-strings table by index),
# columns B-H hold the numeric stats for that country/row.
#
# This edit refreshes the data to a later snapshot (15:28 -> 16:45). Several
# countries that were inserted/re-sorted in the shared-strings table swap
# places with their neighbours (e.g. Azerbaiyan/Guatemala/Honduras), so for
# those rows the label in column A also changes, in addition to the refreshed
# B-H figures. Rows whose country did not move only get new B-H numbers.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Pais")

$ws.Cells.Item(1, 1).Value = 'Datos actualizados a 25 de Junio de 2020 a las 16:45'
$ws.Cells.Item(4, 2).Value = 2464551
$ws.Cells.Item(4, 3).Value = 1997
$ws.Cells.Item(4, 4).Value = 1040669
$ws.Cells.Item(4, 5).Value = 1299551
$ws.Cells.Item(4, 7).Value = 50
$ws.Cells.Item(4, 8).Value = 124331
$ws.Cells.Item(7, 2).Value = 475770
$ws.Cells.Item(7, 3).Value = 2785
$ws.Cells.Item(7, 4).Value = 273270
$ws.Cells.Item(7, 5).Value = 187562
$ws.Cells.Item(7, 7).Value = 31
$ws.Cells.Item(7, 8).Value = 14938
$ws.Cells.Item(8, 2).Value = 307980
$ws.Cells.Item(8, 3).Value = 1118
$ws.Cells.Item(8, 7).Value = 149
$ws.Cells.Item(8, 8).Value = 43230
$ws.Cells.Item(33, 4).Value = 14788
$ws.Cells.Item(33, 5).Value = 33939
$ws.Cells.Item(33, 7).Value = 8
$ws.Cells.Item(33, 8).Value = 1124
$ws.Cells.Item(39, 2).Value = 39139
$ws.Cells.Item(39, 3).Value = 2437
$ws.Cells.Item(39, 4).Value = 18051
$ws.Cells.Item(39, 5).Value = 19651
$ws.Cells.Item(39, 7).Value = 107
$ws.Cells.Item(39, 8).Value = 1437
$ws.Cells.Item(57, 2).Value = 15453
$ws.Cells.Item(57, 3).Value = 375
$ws.Cells.Item(57, 5).Value = 6352
$ws.Cells.Item(59, 1).Value = 'Azerbaiyan'
$ws.Cells.Item(59, 2).Value = 14852
$ws.Cells.Item(59, 3).Value = 547
$ws.Cells.Item(59, 4).Value = 8059
$ws.Cells.Item(59, 5).Value = 6613
$ws.Cells.Item(59, 7).Value = 6
$ws.Cells.Item(59, 8).Value = 180
$ws.Cells.Item(60, 1).Value = 'Guatemala'
$ws.Cells.Item(60, 2).Value = 14819
$ws.Cells.Item(60, 3).Value = 279
$ws.Cells.Item(60, 4).Value = 2930
$ws.Cells.Item(60, 5).Value = 11288
$ws.Cells.Item(60, 7).Value = 19
$ws.Cells.Item(60, 8).Value = 601
$ws.Cells.Item(61, 1).Value = 'Honduras'
$ws.Cells.Item(61, 2).Value = 14571
$ws.Cells.Item(61, 3).Value = 628
$ws.Cells.Item(61, 4).Value = 1546
$ws.Cells.Item(61, 5).Value = 12608
$ws.Cells.Item(61, 7).Value = 12
$ws.Cells.Item(61, 8).Value = 417
$ws.Cells.Item(82, 1).Value = 'Kenia'
$ws.Cells.Item(82, 2).Value = 5384
$ws.Cells.Item(82, 3).Value = 178
$ws.Cells.Item(82, 4).Value = 1857
$ws.Cells.Item(82, 5).Value = 3395
$ws.Cells.Item(82, 7).Value = 2
$ws.Cells.Item(82, 8).Value = 132
$ws.Cells.Item(83, 1).Value = 'El Salvador'
$ws.Cells.Item(83, 2).Value = 5336
$ws.Cells.Item(83, 3).Value = 186
$ws.Cells.Item(83, 4).Value = 3116
$ws.Cells.Item(83, 5).Value = 2094
$ws.Cells.Item(83, 7).Value = 7
$ws.Cells.Item(83, 8).Value = 126
$ws.Cells.Item(84, 1).Value = 'Etiopia'
$ws.Cells.Item(84, 2).Value = 5175
$ws.Cells.Item(84, 3).Value = 141
$ws.Cells.Item(84, 4).Value = 1544
$ws.Cells.Item(84, 5).Value = 3550
$ws.Cells.Item(84, 7).Value = 3
$ws.Cells.Item(84, 8).Value = 81
$ws.Cells.Item(85, 1).Value = 'Guinea'
$ws.Cells.Item(85, 2).Value = 5174
$ws.Cells.Item(85, 4).Value = 3861
$ws.Cells.Item(85, 5).Value = 1284
$ws.Cells.Item(85, 8).Value = 29
$ws.Cells.Item(87, 2).Value = 4635
$ws.Cells.Item(87, 3).Value = 5
$ws.Cells.Item(87, 4).Value = 4269
$ws.Cells.Item(87, 5).Value = 314
$ws.Cells.Item(101, 1).Value = 'Mayotte'
$ws.Cells.Item(101, 2).Value = 2508
$ws.Cells.Item(101, 3).Value = 41
$ws.Cells.Item(101, 4).Value = 2218
$ws.Cells.Item(101, 5).Value = 258
$ws.Cells.Item(101, 8).Value = 32
$ws.Cells.Item(102, 1).Value = 'Croacia'
$ws.Cells.Item(102, 2).Value = 2483
$ws.Cells.Item(102, 3).Value = 95
$ws.Cells.Item(102, 4).Value = 2149
$ws.Cells.Item(102, 5).Value = 227
$ws.Cells.Item(102, 8).Value = 107
$ws.Cells.Item(103, 2).Value = 2321
$ws.Cells.Item(103, 3).Value = 2
$ws.Cells.Item(103, 4).Value = 2171
$ws.Cells.Item(103, 5).Value = 65
$ws.Cells.Item(107, 2).Value = 2010
$ws.Cells.Item(107, 3).Value = 9
$ws.Cells.Item(107, 5).Value = 397
$ws.Cells.Item(111, 1).Value = 'Islandia'
$ws.Cells.Item(111, 2).Value = 1830
$ws.Cells.Item(111, 3).Value = 6
$ws.Cells.Item(111, 4).Value = 1811
$ws.Cells.Item(111, 5).Value = 9
$ws.Cells.Item(111, 8).Value = 10
$ws.Cells.Item(112, 1).Value = 'Madagascar'
$ws.Cells.Item(112, 2).Value = 1829
$ws.Cells.Item(112, 3).Value = 42
$ws.Cells.Item(112, 4).Value = 823
$ws.Cells.Item(112, 5).Value = 990
$ws.Cells.Item(112, 8).Value = 16
$ws.Cells.Item(115, 2).Value = 1662
$ws.Cells.Item(115, 3).Value = 18
$ws.Cells.Item(115, 4).Value = 1144
$ws.Cells.Item(115, 5).Value = 485
$ws.Cells.Item(125, 2).Value = 1162
$ws.Cells.Item(125, 3).Value = 2
$ws.Cells.Item(125, 5).Value = 89
$ws.Cells.Item(147, 1).Value = 'Liberia'
$ws.Cells.Item(147, 2).Value = 681
$ws.Cells.Item(147, 3).Value = 19
$ws.Cells.Item(147, 4).Value = 278
$ws.Cells.Item(147, 5).Value = 369
$ws.Cells.Item(147, 8).Value = 34
$ws.Cells.Item(148, 1).Value = 'Jamaica'
$ws.Cells.Item(148, 2).Value = 678
$ws.Cells.Item(148, 3).Value = 8
$ws.Cells.Item(148, 4).Value = 521
$ws.Cells.Item(148, 5).Value = 147
$ws.Cells.Item(148, 8).Value = 10
$ws.Cells.Item(149, 1).Value = 'Libia'
$ws.Cells.Item(149, 2).Value = 670
$ws.Cells.Item(149, 3).Value = 0
$ws.Cells.Item(149, 4).Value = 138
$ws.Cells.Item(149, 5).Value = 514
$ws.Cells.Item(149, 8).Value = 18
$ws.Cells.Item(150, 1).Value = 'Malta'
$ws.Cells.Item(150, 2).Value = 668
$ws.Cells.Item(150, 3).Value = 3
$ws.Cells.Item(150, 4).Value = 627
$ws.Cells.Item(150, 5).Value = 32
$ws.Cells.Item(150, 8).Value = 9
$ws.Cells.Item(202, 1).Value = 'Dominica'
$ws.Cells.Item(203, 1).Value = 'Fiyi'
$ws.Cells.Item(208, 1).Value = 'Islas Malvinas'
$ws.Cells.Item(209, 1).Value = 'Groenlandia'
$ws.Cells.Item(211, 1).Value = 'Montserrat'
$ws.Cells.Item(211, 4).Value = 10
$ws.Cells.Item(211, 8).Value = 1
$ws.Cells.Item(212, 1).Value = 'Seychelles'
$ws.Cells.Item(212, 4).Value = 11
$ws.Cells.Item(212, 8).Value = 0
